# "Support formula, and dynamic column deep link"
#
# DataSheet:
#   - B4/C4 swap values (25/32 -> 32/21)
#   - new row 7 "TestNameRange": spills the zTurnRange / zTurnFreeLancer
#     named ranges cell-by-cell (implicit intersection)
#   - new row 8 "TestForumla": exercises AVERAGE / MULTIPLY / a shared
#     range formula / a nested IF / a partially-blank range pull
# FormulaSheet:
#   - G2 goes from a broken formula (#ERROR!) to a clean MULTIPLY() call
#   - new C3 with a nested IF()

$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("DataSheet")
$fs = $wb.Worksheets.Item("FormulaSheet")

# ---------------------------------------------------------------------
# DataSheet!B4:C4 - swap the NumOfFreeLancer values for 2008/2009
# ---------------------------------------------------------------------
$ds.Range("B4").Value = 32.0
$ds.Range("C4").Value = 21.0

# ---------------------------------------------------------------------
# Row 7/8 formatting first (copy the existing row 6 look):
#   col A -> same style as A6 ("label" style)
#   cols B.. -> same style as B6 (the shared-formula / formula style)
# ---------------------------------------------------------------------
$ds.Range("A6").Copy()
$ds.Range("A7:A8").PasteSpecial(-4122)

$ds.Range("B6").Copy()
$ds.Range("B7:G7").PasteSpecial(-4122)
$ds.Range("B8:H8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 7 - TestNameRange: one formula per cell (NOT filled as one block)
# so each cell keeps its own independent "=zTurnRange" / "=zTurnFreeLancer"
# formula instead of collapsing into a shared-formula group.
# ---------------------------------------------------------------------
$ds.Range("A7").Value = "TestNameRange"

$ds.Range("B7").Formula = "=zTurnRange"
$ds.Range("C7").Formula = "=zTurnRange"
$ds.Range("D7").Formula = "=zTurnRange"
$ds.Range("E7").Formula = "=zTurnRange"
$ds.Range("F7").Formula = "=zTurnFreeLancer"
$ds.Range("G7").Formula = "=zTurnFreeLancer"

# ---------------------------------------------------------------------
# Row 8 - TestForumla
# ---------------------------------------------------------------------
$ds.Range("A8").Value = "TestForumla"

$ds.Range("B8").Formula = "=AVERAGE(zTurnFreeLancer)"
$ds.Range("C8").Formula = "=MULTIPLY(zTurnRange, zTurnBaseCash)"

# D8:F8 filled together on purpose -> single shared formula, like the
# existing B6:F6 shared formula above it.
$ds.Range("D8:F8").Formula = "=B2:E2"

$ds.Range("G8").Formula = "=If(B5 >= G5, B7, G7)"
$ds.Range("H8").Formula = "=F2:I2"

# ---------------------------------------------------------------------
# FormulaSheet!G2 - was a typo'd formula producing #ERROR!, now a clean
# MULTIPLY call (keeps the existing s="6" style).
# ---------------------------------------------------------------------
$fs.Range("G2").Formula = "=MULTIPLY(zTurnGrowth, zTurnRange)"

# ---------------------------------------------------------------------
# FormulaSheet!C3 - new cell, same style as B3, nested IF().
# ---------------------------------------------------------------------
$fs.Range("B3").Copy()
$fs.Range("C3").PasteSpecial(-4122)
$fs.Range("C3").Formula = "=IF(0>=2, 5, 10)"

$ds.Range("A1").Select()
